$wb = $excel.ActiveWorkbook

# --- Sheet "Statistics": remove the three sample data rows (rows 2-4), ---
# --- keeping only the header row.                                      ---
$ws1 = $wb.Worksheets.Item("Statistics")
$ws1.Range("A2:C4").EntireRow.Delete()

# --- Sheet "Accidents": append the simulated accident rows (rows 2-19) ---
$ws2 = $wb.Worksheets.Item("Accidents")

$accidentRows = @(
    @("2024-07-27 11:30:28", "Truck and Car", "57.40 and 90.70"),
    @("2024-07-27 11:30:28", "Car and Car", "114.81 and 95.71"),
    @("2024-07-27 11:30:28", "Car and Car", "104.73 and 0.00"),
    @("2024-07-27 11:30:28", "Car and Car", "103.94 and 95.46"),
    @("2024-07-27 11:30:29", "Car and Car", "85.61 and 0.00"),
    @("2024-07-27 11:30:29", "Car and Car", "102.74 and 0.00"),
    @("2024-07-27 11:30:29", "Car and Car", "89.08 and 0.00"),
    @("2024-07-27 11:30:30", "Car and Car", "103.45 and 0.00"),
    @("2024-07-27 11:30:30", "Car and Car", "107.26 and 0.00"),
    @("2024-07-27 11:30:30", "Car and Car", "66.65 and 0.00"),
    @("2024-07-27 11:30:30", "Car and Car", "107.40 and 0.00"),
    @("2024-07-27 11:30:30", "Car and Car", "80.07 and 0.00"),
    @("2024-07-27 11:30:30", "Car and Car", "99.26 and 0.00"),
    @("2024-07-27 11:30:31", "Car and Car", "100.17 and 0.00"),
    @("2024-07-27 11:30:31", "Truck and Car", "98.51 and 0.00"),
    @("2024-07-27 11:30:31", "Truck and Car", "66.08 and 0.00"),
    @("2024-07-27 11:30:32", "Car and Truck", "64.28 and 0.00"),
    @("2024-07-27 11:30:32", "Car and Truck", "92.91 and 0.00")
)

$r = 2
foreach ($row in $accidentRows) {
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}
